$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (pushes the existing rows 8-20 down to 9-21),
# inheriting the formatting of the surrounding rows.
$ws.Rows.Item(8).Insert()

# Row 7 used to describe "t_find_originals_and_duplicates" (a module that is
# gone this week); it now describes the newly introduced
# "t_delete_empty_variales" module instead.
$ws.Range("A7").Value = "t_delete_empty_variales"
$ws.Range("B7").Value = "not available"
$ws.Range("C7").Value = "not vailable"
$ws.Range("D7").Value = "introduced"

# The freshly inserted row 8 reports the same status as the top-level
# "r_toolbox" package row (introduced / kept / changed), with no module name.
$ws.Range("B8").Value = "introduced"
$ws.Range("C8").Value = "kept"
$ws.Range("D8").Value = "changed"

# Reflect the last-used cell of the bottom-right pane at save time.
$ws.Range("D8").Select()
